$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75; this pushes the existing rows 75..197
# down to 76..198 (matching the diff's "everything shifts down by one,
# with the old last row (197) now living at row 198" pattern).
$ws.Rows.Item(75).Insert()

# Populate the freshly inserted row 75 with the new record's data.
$ws.Range("A75").Value = 3
$ws.Range("B75").Value = 'Femacal de La Calera'
$ws.Range("C75").Value = 'Coquimbo'
$ws.Range("D75").Value = 44477
$ws.Range("E75").Value = 5
$ws.Range("F75").Value = 100112012
$ws.Range("G75").Value = 'Espinaca'
$ws.Range("H75").Value = 'Sin especificar'
$ws.Range("I75").Value = 'Primera'
$ws.Range("J75").Value = 160
$ws.Range("K75").Value = 2500
$ws.Range("L75").Value = 2500
$ws.Range("M75").Value = 2500
$ws.Range("N75").Value = '$/docena de atados (3 kilos)'
$ws.Range("O75").Value = 'Provincia de Quillota'
$ws.Range("P75").Value = 833
$ws.Range("Q75").Value = 3
$ws.Range("R75").Value = 'Hortaliza'
